$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.432.85"
$ws.Range("E2").Value = "  -0.10%  "

$ws.Range("D3").Value = "2.227.63"
$ws.Range("E3").Value = "  -1.03%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'298.83"
$ws.Range("E5").Value = "  -3.19%  "

$ws.Range("D6").Value = "'90.28"
$ws.Range("E6").Value = "  -5.08%  "

$ws.Range("D7").Value = "'0.562"
$ws.Range("E7").Value = "  -1.77%  "

$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("D9").Value = "'0.493"
$ws.Range("E9").Value = "  -6.54%  "

$ws.Range("D10").Value = "'33.07"
$ws.Range("E10").Value = "  -6.53%  "

$ws.Range("D11").Value = "'0.0784"
$ws.Range("E11").Value = "  -3.31%  "

$ws.Range("D12").Value = "'6.96"
$ws.Range("E12").Value = "  -4.39%  "

$ws.Range("E13").Value = "  -0.55%  "

$ws.Range("D14").Value = "2.564.70"
$ws.Range("E14").Value = "  -1.10%  "

$ws.Range("D15").Value = "2.283.01"
$ws.Range("E15").Value = "  -1.35%  "

$ws.Range("D16").Value = "'13.47"
$ws.Range("E16").Value = "  -1.59%  "

$ws.Range("D17").Value = "'0.782"
$ws.Range("E17").Value = "  -7.00%  "

$ws.Range("D18").Value = "44.204.35"
$ws.Range("E18").Value = "  -0.03%  "

$ws.Range("D19").Value = "0.0₃0912"
$ws.Range("E19").Value = "  -5.78%  "

$ws.Range("D20").Value = "'5.93"
$ws.Range("E20").Value = "  -7.89%  "

$ws.Range("D21").Value = "'11.01"
$ws.Range("E21").Value = "  -10.10%  "

$ws.Range("D22").Value = "'64.63"
$ws.Range("E22").Value = "  -2.18%  "

$ws.Range("D23").Value = "'240.53"
$ws.Range("E23").Value = "  +1.33%  "

$ws.Range("D24").Value = "'2.80"
$ws.Range("E24").Value = "  -6.81%  "

$ws.Range("E25").Value = "  +1.04%  "

$ws.Range("E26").Value = "  -7.62%  "

$ws.Range("D27").Value = "'2.25"
$ws.Range("E27").Value = "  +0.91%  "

$ws.Range("D28").Value = "'38.72"
$ws.Range("E28").Value = "  +0.89%  "

$ws.Range("D29").Value = "'9.39"
$ws.Range("E29").Value = "  -4.92%  "

$ws.Range("D30").Value = "'19.49"
$ws.Range("E30").Value = "  -3.12%  "

$ws.Range("D31").Value = "'148.59"
$ws.Range("E31").Value = "  -2.91%  "

$ws.Range("D32").Value = "'5.40"
$ws.Range("E32").Value = "  -10.11%  "

$ws.Range("D33").Value = "'2.54"
$ws.Range("E33").Value = "  -3.55%  "

$ws.Range("D34").Value = "'0.0753"
$ws.Range("E34").Value = "  -5.90%  "

$ws.Range("E35").Value = "  -3.89%  "

$ws.Range("E36").Value = "  -10.24%  "

$ws.Range("E37").Value = "  -6.25%  "

$ws.Range("D38").Value = "'1.68"
$ws.Range("E38").Value = "  -7.06%  "

$ws.Range("D39").Value = "'0.0306"
$ws.Range("E39").Value = "  +0.98%  "

$ws.Range("D40").Value = "'3.20"
$ws.Range("E40").Value = "  -7.42%  "

$ws.Range("D41").Value = "'3.54"
$ws.Range("E41").Value = "  -7.64%  "

$ws.Range("D42").Value = "'13.16"
$ws.Range("E42").Value = "  -10.36%  "

$ws.Range("E43").Value = "  -0.24%  "

$ws.Range("D44").Value = "1.820.37"
$ws.Range("E44").Value = "  +3.76%  "

$ws.Range("D45").Value = "'1.78"
$ws.Range("E45").Value = "  +10.95%  "

$ws.Range("E46").Value = "  -7.71%  "

$ws.Range("B47").Value = "BitcoinSV"
$ws.Range("C47").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D47").Value = "'74.92"
$ws.Range("E47").Value = "  -7.57%  "

$ws.Range("B48").Value = "HuobiToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D48").Value = "'2.83"
$ws.Range("E48").Value = "  +13.73%  "

$ws.Range("D49").Value = "'94.53"
$ws.Range("E49").Value = "  -5.70%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'14.12"
$ws.Range("E50").Value = "  +7.89%  "

$ws.Range("B51").Value = "ordi"
$ws.Range("C51").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D51").Value = "'67.07"
$ws.Range("E51").Value = "  -6.31%  "
